$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Datos" to "lote"
$ws.Name = "lote"

# Remove the bold/white-on-blue header styling that used to highlight row 1
$ws.Range("A1:E1").ClearFormats()

# Rewrite the header row: lower-case labels and a new "finca" column
# inserted after "nombre", pushing descripcion/criterio/comentario over
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "finca"
$ws.Range("D1").Value = "descripcion"
$ws.Range("E1").Value = "criterio"
$ws.Range("F1").Value = "comentario"
